# A new weekly price record was inserted into the "Hortaliza, Macroferia
# Regional de Talca - Papa" data table. All rows from 254 downward shift
# by one; row 254 becomes the newly-reported record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 254 - this pushes the existing
# rows 254..312 down to 255..313 (old row 312's data lands on new row 313
# automatically, exactly matching the target state).
$ws.Rows("254:254").Insert()

# Populate the newly inserted row 254 with the new record.
$ws.Cells.Item(254, 1).Value = 5
$ws.Cells.Item(254, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(254, 3).Value = "Maule"
$ws.Cells.Item(254, 4).Value = 44508
$ws.Cells.Item(254, 5).Value = 7
$ws.Cells.Item(254, 6).Value = 100114001
$ws.Cells.Item(254, 7).Value = "Papa"
$ws.Cells.Item(254, 8).Value = "Rodeo"
$ws.Cells.Item(254, 9).Value = "1a nueva(o)"
$ws.Cells.Item(254, 10).Value = 1600
$ws.Cells.Item(254, 11).Value = 10000
$ws.Cells.Item(254, 12).Value = 10000
$ws.Cells.Item(254, 13).Value = 10000
$ws.Cells.Item(254, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(254, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(254, 16).Value = 400
$ws.Cells.Item(254, 17).Value = 25
$ws.Cells.Item(254, 18).Value = "Hortaliza"
